$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The sheet previously held a small "property / value" table (A1:B5, with
# some boolean TRUE cells). It is being replaced by a "spice / superclass"
# table listing 23 spices & herbs together with the cuisine group(s) they
# belong to (mexican / indian / asian / pass). One spice (cumin) belongs to
# two groups, so it also gets a value in column C.
# ---------------------------------------------------------------------------

# Clear out everything that is currently on the sheet before laying down the
# new table, so no stale cells (e.g. the old boolean TRUE values) survive.
$ws.Cells.Clear()

$data = @(
    @("cardamom", "asian"),
    @("allspice", "pass"),
    @("garlic", "mexican"),
    @("cayenne pepper", "indian"),
    @("ginger", "asian"),
    @("coriander", "mexican"),
    @("ras el hanout", "pass"),
    @("turmeric", "asian"),
    @("curry", "indian"),
    @("cinnamon", "mexican"),
    @("herbes de province", "pass"),
    @("cumin", "mexican", "asian"),
    @("rosemary", "pass"),
    @("bay leaves", "indian"),
    @("basil", "asian"),
    @("chili powder", "mexican"),
    @("thyme", "pass"),
    @("nutmeg", "indian"),
    @("cajun seasoning", "pass"),
    @("za'atar", "pass"),
    @("garam masala", "indian"),
    @("oregano", "mexican"),
    @("paprika", "pass")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $i + 1
    $entry = $data[$i]

    $ws.Cells.Item($rowNum, 1).Value = $entry[0]
    $ws.Cells.Item($rowNum, 2).Value = $entry[1]

    if ($entry.Count -gt 2) {
        $ws.Cells.Item($rowNum, 3).Value = $entry[2]
    }
}

# Column A should auto-fit its longest entry ("herbes de province").
$ws.Columns.Item(1).ColumnWidth = 15.6

# Leave the selection just below the table, matching where a user would end
# up after typing the last row and pressing Enter/Down.
$ws.Range("B24").Select() | Out-Null
